# Auto-generated edit script: updates crypto price/volume data
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "82.013.49"
$ws.Range("E2").Value = "  +3.11%  "
$ws.Range("D3").Value = "3.191.49"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'215.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.12%  "
$ws.Range("D6").Value = "'626.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +20.58%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").Value = "3.189.92"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'0.592"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'0.0000259"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.76%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'5.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.17%  "
$ws.Range("D15").Value = "3.781.57"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "81.843.68"
$ws.Range("D18").Value = "3.191.72"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +6.46%  "
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "'435.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "'8.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.69%  "
$ws.Range("D23").Value = "'5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'7.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.52%  "
$ws.Range("D25").Value = "'5.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.94%  "
$ws.Range("D26").Value = "3.350.81"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'76.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "'10.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +4.14%  "
$ws.Range("D31").Value = "'586.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.78%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'9.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  +7.14%  "
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("E37").Value = "  +15.60%  "
$ws.Range("D38").Value = "'22.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'6.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.45%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  +13.78%  "
$ws.Range("D43").Value = "'3.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +22.80%  "
$ws.Range("D44").Value = "'20.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.74%  "
$ws.Range("D45").Value = "'160.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'188.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "'44.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "'0.775"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.10%  "
$ws.Range("D51").Value = "'26.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "
